$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.46865769738307961
$ws.Range("A2").Value = -0.0099999996771309441
$ws.Range("A3").Value = -0.008999999668443337
$ws.Range("A4").Value = -0.01199999990735634
$ws.Range("A5").Value = -0.0059999996748691942
$ws.Range("A6").Value = -0.0059999996673241185
$ws.Range("A7").Value = -0.019999999603568241
$ws.Range("A8").Value = -0.019999999602110741
$ws.Range("A9").Value = -0.005999999664513922
$ws.Range("A10").Value = -0.0059999996643753661
$ws.Range("A11").Value = 0.054667101718045785
$ws.Range("A12").Value = -0.0059999996630351049
$ws.Range("A13").Value = -0.0059999996577113635
$ws.Range("A14").Value = -0.011999999628996783
$ws.Range("A15").Value = -0.0059999996547466239
$ws.Range("A16").Value = -0.0059999996535355926
$ws.Range("A17").Value = -0.0059999996519524146
$ws.Range("A18").Value = -0.0089999996380356606
$ws.Range("A19").Value = -0.083388802885858482
$ws.Range("A20").Value = -0.0089999996773251212
$ws.Range("A21").Value = -0.0089999996768446167
$ws.Range("A22").Value = -0.0089999996764866808
$ws.Range("A23").Value = -0.06218577027784189
$ws.Range("A24").Value = -0.04199999950367328
$ws.Range("A25").Value = -0.087777373258987268
$ws.Range("A26").Value = -0.0059999996665389688
$ws.Range("A27").Value = -0.0059999996653621324
$ws.Range("A28").Value = -0.0059999996611006523
$ws.Range("A29").Value = -0.011999999631409963
$ws.Range("A30").Value = -0.019999999594255691
$ws.Range("A31").Value = 0.01745964040600434
$ws.Range("A32").Value = -0.020999999588400264
$ws.Range("A33").Value = -0.0059999996557786872
